# "combine genomic traits and experiments into one big table"
#
# The "ROS detoxification" experimental-rate columns (D:I, rows 3-10) were
# entered with the wrong unit scale (x 1e-6 too small) relative to the
# "proteolytic" sheet's rates; correct them by rescaling by 1,000,000 and
# bring the "sd/COV" column's number format in line with its neighbours.
# Extra (currently blank) rows are pre-formatted below the table, ready to
# receive more combined assay rows, and the active sheet/selection is moved
# onto the corrected table.

$wb = $excel.ActiveWorkbook
$wsRos = $wb.Worksheets.Item("ROS detoxification")
$wsProt = $wb.Worksheets.Item("proteolytic")

# --- Fix the unit scale on the raw + derived rate columns (D:I, rows 3-10) ---
$rates = $wsRos.Range("D3:I10")
foreach ($cell in $rates) {
    $cell.Value() = $cell.Value() * 1000000
}

# Column I (COV) was left as "General" on several rows; match the scientific
# notation used by the rest of the table now that the values are corrected.
$wsRos.Range("I3:I10").NumberFormat() = "0.00E+00"

# --- Pre-format blank rows under the table for the combined data to come ---
$wsRos.Range("D13:I22").NumberFormat() = "0.00E+00"

# --- Move focus to the "ROS detoxification" sheet / corrected table ---
$wsRos.Activate()
$wsRos.Range("J3").Select()
